$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing the existing row 7 (and everything
# below it) down by one. This is a weekly data refresh: a new day's record is
# added at the top of the dated series and every older record shifts down one
# slot (the oldest record that falls off the bottom is re-appended as the new
# last row, row 142).
$ws.Rows(7).Insert()

# Populate the newly-inserted row 7 with the new record's data. The
# non-varying descriptive columns match every other row in this sheet.
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = 44691
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 100112012
$ws.Range("G7").Value = "Espinaca"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 65
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = 9000
$ws.Range("N7").Value = "$/docena de atados"
$ws.Range("O7").Value = "Región de La Araucanía"
$ws.Range("P7").Value = 3000
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = "Hortaliza"
